$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 350, shifting existing rows 350:412 down to 351:413
$ws.Rows.Item(350).Insert()

# Populate the newly inserted row 350 with the new data record
$ws.Cells.Item(350, 1).Value = 4
$ws.Cells.Item(350, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(350, 3).Value = "Los Lagos"
$ws.Cells.Item(350, 4).Value = "10/12/2023"
$ws.Cells.Item(350, 5).Value = 10
$ws.Cells.Item(350, 6).Value = 100112039
$ws.Cells.Item(350, 7).Value = "Ciboulette"
$ws.Cells.Item(350, 8).Value = "Sin especificar"
$ws.Cells.Item(350, 9).Value = "Primera"
$ws.Cells.Item(350, 10).Value = 80
$ws.Cells.Item(350, 11).Value = 3500
$ws.Cells.Item(350, 12).Value = 3500
$ws.Cells.Item(350, 13).Value = 3500
$ws.Cells.Item(350, 14).Value = "`$/docena de atados"
$ws.Cells.Item(350, 15).Value = "Región Metropolitana"
$ws.Cells.Item(350, 16).Value = 1167
$ws.Cells.Item(350, 17).Value = 3
$ws.Cells.Item(350, 18).Value = "Hortaliza"
